$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "292.63"
Set-TextValue $ws.Range("E2") "-6.70%"
Set-TextValue $ws.Range("G2") "7"
Set-TextValue $ws.Range("D3") "40.44"
Set-TextValue $ws.Range("E3") "-1.05%"
Set-TextValue $ws.Range("G3") "7"
Set-TextValue $ws.Range("D4") "5.044"
Set-TextValue $ws.Range("E4") "-1.95%"
Set-TextValue $ws.Range("G4") "7"
Set-TextValue $ws.Range("E5") "-3.62%"
Set-TextValue $ws.Range("G5") "7"
Set-TextValue $ws.Range("D6") "1.532"
Set-TextValue $ws.Range("E6") "-8.85%"
Set-TextValue $ws.Range("G6") "7"
Set-TextValue $ws.Range("D7") "0.9301"
Set-TextValue $ws.Range("E7") "-0.08%"
Set-TextValue $ws.Range("G7") "7"
Set-TextValue $ws.Range("D8") "2.358"
Set-TextValue $ws.Range("E8") "-2.72%"
Set-TextValue $ws.Range("G8") "7"
Set-TextValue $ws.Range("D9") "0.1174"
Set-TextValue $ws.Range("E9") "-2.03%"
Set-TextValue $ws.Range("G9") "7"
Set-TextValue $ws.Range("D10") "0.1741"
Set-TextValue $ws.Range("E10") "-4.05%"
Set-TextValue $ws.Range("G10") "7"
Set-TextValue $ws.Range("D11") "0.04347"
Set-TextValue $ws.Range("E11") "5.02%"
Set-TextValue $ws.Range("G11") "7"
Set-TextValue $ws.Range("D12") "0.08705"
Set-TextValue $ws.Range("E12") "-3.65%"
Set-TextValue $ws.Range("G12") "7"
Set-TextValue $ws.Range("E13") "0.00%"
Set-TextValue $ws.Range("G13") "7"
Set-TextValue $ws.Range("D14") "0.001273"
Set-TextValue $ws.Range("E14") "-1.58%"
Set-TextValue $ws.Range("G14") "7"
Set-TextValue $ws.Range("D15") "0.005956"
Set-TextValue $ws.Range("E15") "2.01%"
Set-TextValue $ws.Range("G15") "7"
Set-TextValue $ws.Range("E16") "-0.02%"
Set-TextValue $ws.Range("G16") "7"
Set-TextValue $ws.Range("D17") "4.283"
Set-TextValue $ws.Range("E17") "-0.96%"
Set-TextValue $ws.Range("G17") "7"
Set-TextValue $ws.Range("D18") "0.3290"
Set-TextValue $ws.Range("E18") "-1.65%"
Set-TextValue $ws.Range("G18") "7"
Set-TextValue $ws.Range("D19") "7.971"
Set-TextValue $ws.Range("E19") "4.76%"
Set-TextValue $ws.Range("G19") "7"
Set-TextValue $ws.Range("E20") "3.58%"
Set-TextValue $ws.Range("G20") "7"
Set-TextValue $ws.Range("D21") "0.2743"
Set-TextValue $ws.Range("E21") "-3.47%"
Set-TextValue $ws.Range("G21") "7"
Set-TextValue $ws.Range("D22") "0.03931"
Set-TextValue $ws.Range("E22") "-1.08%"
Set-TextValue $ws.Range("G22") "7"
Set-TextValue $ws.Range("D23") "0.001262"
Set-TextValue $ws.Range("E23") "-1.36%"
Set-TextValue $ws.Range("G23") "7"
Set-TextValue $ws.Range("D24") "0.003786"
Set-TextValue $ws.Range("E24") "-7.38%"
Set-TextValue $ws.Range("G24") "7"
Set-TextValue $ws.Range("E25") "-5.32%"
Set-TextValue $ws.Range("G25") "7"
Set-TextValue $ws.Range("D26") "0.0003725"
Set-TextValue $ws.Range("G26") "7"
Set-TextValue $ws.Range("G27") "7"
Set-TextValue $ws.Range("G28") "7"
Set-TextValue $ws.Range("G29") "7"
Set-TextValue $ws.Range("G30") "7"
Set-TextValue $ws.Range("G31") "7"
Set-TextValue $ws.Range("G32") "7"
Set-TextValue $ws.Range("G33") "7"
Set-TextValue $ws.Range("G34") "7"
Set-TextValue $ws.Range("G35") "7"
Set-TextValue $ws.Range("G36") "7"
Set-TextValue $ws.Range("G37") "7"
Set-TextValue $ws.Range("D38") "0.02289"
Set-TextValue $ws.Range("E38") "-4.94%"
Set-TextValue $ws.Range("G38") "7"
Set-TextValue $ws.Range("D39") "0.05059"
Set-TextValue $ws.Range("E39") "-1.76%"
Set-TextValue $ws.Range("G39") "7"
Set-TextValue $ws.Range("D40") "0.006325"
Set-TextValue $ws.Range("E40") "91.30%"
Set-TextValue $ws.Range("G40") "7"
Set-TextValue $ws.Range("D41") "0.007818"
Set-TextValue $ws.Range("E41") "0.92%"
Set-TextValue $ws.Range("G41") "7"
Set-TextValue $ws.Range("D42") "0.1287"
Set-TextValue $ws.Range("E42") "-1.18%"
Set-TextValue $ws.Range("G42") "7"
Set-TextValue $ws.Range("D43") "0.007378"
Set-TextValue $ws.Range("E43") "-2.88%"
Set-TextValue $ws.Range("G43") "7"
Set-TextValue $ws.Range("D44") "0.008259"
Set-TextValue $ws.Range("E44") "-2.77%"
Set-TextValue $ws.Range("G44") "7"
Set-TextValue $ws.Range("D45") "0.2924"
Set-TextValue $ws.Range("E45") "-13.66%"
Set-TextValue $ws.Range("G45") "7"
Set-TextValue $ws.Range("D46") "0.00006286"
Set-TextValue $ws.Range("E46") "-4.66%"
Set-TextValue $ws.Range("G46") "7"
Set-TextValue $ws.Range("E47") "-0.11%"
Set-TextValue $ws.Range("G47") "7"
Set-TextValue $ws.Range("D48") "0.03261"
Set-TextValue $ws.Range("E48") "-87.86%"
Set-TextValue $ws.Range("G48") "7"
Set-TextValue $ws.Range("E49") "-0.11%"
Set-TextValue $ws.Range("G49") "7"
Set-TextValue $ws.Range("D50") "0.0002001"
Set-TextValue $ws.Range("E50") "-0.11%"
Set-TextValue $ws.Range("G50") "7"
Set-TextValue $ws.Range("G51") "7"
